# Update cryptocurrency price/volume table to the latest scraped values.
# Also fixes two swapped row pairs (Chainlink/TRON at rows 14-15, and
# Cronos/PaxDollar at rows 47-48) where the source site's ranking order changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.228.80"
$ws.Range("E2").Value = "  -3.47%  "

# Row 3
$ws.Range("D3").Value = "1.806.38"
$ws.Range("E3").Value = "  -3.78%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.56"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4210"
$ws.Range("E7").Value = "  -2.51%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3548"
$ws.Range("E8").Value = "  -4.24%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07127"
$ws.Range("E9").Value = "  -4.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8450"
$ws.Range("E10").Value = "  -4.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.14"
$ws.Range("E11").Value = "  -5.06%  "

# Row 12
$ws.Range("D12").Value = "1.792.96"
$ws.Range("E12").Value = "  -3.81%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.325"
$ws.Range("E13").Value = "  -2.98%  "

# Row 14
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.06919"
$ws.Range("E14").Value = "  -1.08%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.354"
$ws.Range("E15").Value = "  -4.39%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.07"
$ws.Range("E17").Value = "  -0.20%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008791"
$ws.Range("E18").Value = "  -3.84%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.26%  "

# Row 20
$ws.Range("E20").Value = "  -3.49%  "

# Row 21
$ws.Range("D21").Value = "27.535.90"
$ws.Range("E21").Value = "  -2.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.074"
$ws.Range("E22").Value = "  -0.49%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  -1.21%  "

# Row 24
$ws.Range("D24").Value = "2.090.63"
$ws.Range("E24").Value = "  -2.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.962"
$ws.Range("E25").Value = "  -1.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.15"
$ws.Range("E26").Value = "  -0.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.20"
$ws.Range("E27").Value = "  -2.71%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.056"
$ws.Range("E28").Value = "  -7.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.20"
$ws.Range("E29").Value = "  -4.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.716"
$ws.Range("E30").Value = "  -9.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08892"
$ws.Range("E31").Value = "  -1.12%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7401"
$ws.Range("E32").Value = "  -7.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.949"
$ws.Range("E33").Value = "  -1.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.462"
$ws.Range("E34").Value = "  -5.08%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.102"
$ws.Range("E35").Value = "  -6.36%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.004"
$ws.Range("E36").Value = "  +0.08%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.068"
$ws.Range("E37").Value = "  -5.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05213"
$ws.Range("E38").Value = "  -4.92%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01901"
$ws.Range("E39").Value = "  -3.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.760"
$ws.Range("E40").Value = "  -4.64%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1639"
$ws.Range("E41").Value = "  -3.65%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4982"
$ws.Range("E42").Value = "  -3.82%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.286"
$ws.Range("E43").Value = "  -8.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.197"
$ws.Range("E44").Value = "  -4.73%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.32"
$ws.Range("E45").Value = "  -2.46%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.84"
$ws.Range("E46").Value = "  -0.85%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.003"
$ws.Range("E47").Value = "  +0.10%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06407"
$ws.Range("E48").Value = "  -3.02%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4572"
$ws.Range("E49").Value = "  -4.31%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.595"
$ws.Range("E50").Value = "  -4.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.08"
$ws.Range("E51").Value = "  -3.52%  "
